$wb = $excel.ActiveWorkbook

# --- Sheet: Full results ---
$ws1 = $wb.Worksheets.Item("Full results")

# Row 2 (COMPLETE MODEL)
$ws1.Range("H2").Value = 0.589155258229011
$ws1.Range("I2").Value = 0.279370328170232
$ws1.Range("O2").Value = 0.410898045795768

# Row 3 (CONDITIONAL MODEL)
$ws1.Range("F3").Value = 0.603924329216196
$ws1.Range("G3").Value = 0.311741230630201

# Row 4 (NULL MODEL)
$ws1.Range("C4").Value = 0.63758128463405
$ws1.Range("D4").Value = 0.362509198923953
$ws1.Range("E4").Value = 1.000090483558
$ws1.Range("J4").Value = 0.362476400765031
$ws1.Range("K4").Value = 0.31171302572226
$ws1.Range("L4").Value = 0.0147677347384245
$ws1.Range("M4").Value = 0.0484216450307373
$ws1.Range("N4").Value = 0.326480760460685

# --- Sheet: For plotting ---
$ws2 = $wb.Worksheets.Item("For plotting")

# Row 2 (Sibcorr, education)
$ws2.Range("C2").Value = 0.362476400765031
$ws2.Range("D2").Value = 0.333980048443841
$ws2.Range("E2").Value = 0.390972753086221

# Row 3 (IOLIB, education)
$ws2.Range("C3").Value = 0.326480760460685
$ws2.Range("D3").Value = 0.299440473577494
$ws2.Range("E3").Value = 0.353521047343875

# Row 4 (IORAD, education)
$ws2.Range("C4").Value = 0.410898045795768
$ws2.Range("D4").Value = 0.380231140088459
$ws2.Range("E4").Value = 0.441564951503078
